$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.1893333333333333
$ws.Range("C2").Value = 0.5333333333333333
$ws.Range("J2").Value = 0.03733333333333334
$ws.Range("P2").Value = 0.1573333333333333
$ws.Range("S2").Value = 0.08266666666666667
$ws.Range("B3").Value = 0.02403846153846154
$ws.Range("C3").Value = 0.03365384615384615
$ws.Range("J3").Value = 0.04326923076923077
$ws.Range("P3").Value = 0.7211538461538461
$ws.Range("S3").Value = 0.1778846153846154
$ws.Range("J4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2619047619047619
$ws.Range("P5").Value = 0.8333333333333334
$ws.Range("S5").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.09486166007905138
$ws.Range("D6").Value = 0.01976284584980237
$ws.Range("E6").Value = 0.003952569169960474
$ws.Range("F6").Value = 0.09881422924901186
$ws.Range("J6").Value = 0.225296442687747
$ws.Range("O6").Value = 0.03162055335968379
$ws.Range("Q6").Value = 0.150197628458498
$ws.Range("R6").Value = 0.04347826086956522
$ws.Range("S6").Value = 0.3320158102766799
$ws.Range("B7").Value = 0.1483516483516484
$ws.Range("D7").Value = 0.01098901098901099
$ws.Range("E7").Value = 0.01098901098901099
$ws.Range("F7").Value = 0.02197802197802198
$ws.Range("J7").Value = 0.1813186813186813
$ws.Range("O7").Value = 0.03846153846153846
$ws.Range("Q7").Value = 0.1373626373626374
$ws.Range("R7").Value = 0.05494505494505494
$ws.Range("S7").Value = 0.3956043956043956
$ws.Range("B8").Value = 0.1
$ws.Range("D8").Value = 0.0125
$ws.Range("E8").Value = 0.0025
$ws.Range("F8").Value = 0.065
$ws.Range("J8").Value = 0.115
$ws.Range("O8").Value = 0.02
$ws.Range("Q8").Value = 0.23
$ws.Range("R8").Value = 0.06
$ws.Range("S8").Value = 0.395
$ws.Range("B9").Value = 0.09836065573770492
$ws.Range("D9").Value = 0.02459016393442623
$ws.Range("F9").Value = 0.04918032786885246
$ws.Range("J9").Value = 0.1024590163934426
$ws.Range("O9").Value = 0.01229508196721311
$ws.Range("Q9").Value = 0.2213114754098361
$ws.Range("R9").Value = 0.04918032786885246
$ws.Range("S9").Value = 0.4426229508196721
$ws.Range("B10").Value = 0.1268709907341411
$ws.Range("D10").Value = 0.02066999287241625
$ws.Range("E10").Value = 0.00142551674982181
$ws.Range("F10").Value = 0.07697790449037777
$ws.Range("J10").Value = 0.1240199572344975
$ws.Range("O10").Value = 0.02209550962223806
$ws.Range("Q10").Value = 0.2394868139700642
$ws.Range("R10").Value = 0.05773342836778332
$ws.Range("S10").Value = 0.33071988595866
$ws.Range("G11").Value = 0.1414790996784566
$ws.Range("J11").Value = 0.09967845659163987
$ws.Range("K11").Value = 0.1929260450160772
$ws.Range("L11").Value = 0.5594855305466238
$ws.Range("S11").Value = 0.006430868167202572
$ws.Range("G12").Value = 0.7237569060773481
$ws.Range("J12").Value = 0.2154696132596685
$ws.Range("K12").Value = 0.005524861878453038
$ws.Range("L12").Value = 0.02209944751381215
$ws.Range("S12").Value = 0.03314917127071823
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.3928571428571428
$ws.Range("S13").Value = 0.1071428571428571
$ws.Range("F15").Value = 0.01639344262295082
$ws.Range("H15").Value = 0.139344262295082
$ws.Range("I15").Value = 0.09836065573770492
$ws.Range("J15").Value = 0.3442622950819672
$ws.Range("K15").Value = 0.04918032786885246
$ws.Range("M15").Value = 0.00819672131147541
$ws.Range("O15").Value = 0.08196721311475409
$ws.Range("S15").Value = 0.2622950819672131
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.1404255319148936
$ws.Range("I16").Value = 0.1063829787234043
$ws.Range("J16").Value = 0.4340425531914894
$ws.Range("K16").Value = 0.08085106382978724
$ws.Range("M16").Value = 0.01702127659574468
$ws.Range("N16").Value = 0.00425531914893617
$ws.Range("O16").Value = 0.05531914893617021
$ws.Range("S16").Value = 0.1404255319148936
$ws.Range("F17").Value = 0.02402957486136784
$ws.Range("H17").Value = 0.1626617375231054
$ws.Range("I17").Value = 0.1275415896487985
$ws.Range("J17").Value = 0.4343807763401109
$ws.Range("K17").Value = 0.07948243992606285
$ws.Range("M17").Value = 0.007393715341959334
$ws.Range("O17").Value = 0.05545286506469501
$ws.Range("S17").Value = 0.1090573012939002
$ws.Range("F18").Value = 0.0364963503649635
$ws.Range("H18").Value = 0.1751824817518248
$ws.Range("I18").Value = 0.08759124087591241
$ws.Range("J18").Value = 0.4379562043795621
$ws.Range("K18").Value = 0.072992700729927
$ws.Range("M18").Value = 0.0218978102189781
$ws.Range("O18").Value = 0.0583941605839416
$ws.Range("S18").Value = 0.1094890510948905
$ws.Range("F19").Value = 0.02021772939346812
$ws.Range("H19").Value = 0.1734059097978227
$ws.Range("I19").Value = 0.09020217729393468
$ws.Range("J19").Value = 0.3856920684292379
$ws.Range("K19").Value = 0.1228615863141524
$ws.Range("M19").Value = 0.01477449455676516
$ws.Range("N19").Value = 0.001555209953343701
$ws.Range("O19").Value = 0.06920684292379471
$ws.Range("S19").Value = 0.1220839813374806
